$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new citation row (row 21). Set values in the same order the shared
# strings are introduced in the source file: Title (83), Authors (84),
# Publication (85), Pages (86). Publisher ("Elsevier") already exists in
# the shared-string table and is reused.
$ws.Range("B21").Value = "A modular approach for the ultra-scale-down of depth filtration"
$ws.Range("A21").Value = "Aaron Noyes, Jonida Basha, John Frostad, Scott Cook, Doug Millard, Jim Mullin, Daniel LaCasse, Richard S Wright, Benjamin Huffman, Robert Fahrner, Ranga Godavarti, Nigel Titchener-Hooker, Khurram Sunasara, Tarit Mukhopadhyay"
$ws.Range("C21").Value = "Journal of Membrane Science"
$ws.Range("D21").Value = 496
$ws.Range("F21").Value = "199-210"
$ws.Range("G21").Value = 2015
$ws.Range("H21").Value = "Elsevier"

# Move the active selection, matching the author's final cursor position.
$ws.Range("C27").Select()
